# Update NATMI ligand-receptor statistics (Cd14-Itgb1) with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.331241
$ws.Range("H2").Value = 0.9937229999999999
$ws.Range("I2").Value = 0.2019242199214145
$ws.Range("J2").Value = 0.2019242199214145
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 20.21948272345967
$ws.Range("R2").Value = 181.975344511137
$ws.Range("S2").Value = 0.04126550538474153
$ws.Range("T2").Value = 0.04126550538474152
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.331241
$ws.Range("H3").Value = 0.9937229999999999
$ws.Range("I3").Value = 0.2019242199214145
$ws.Range("J3").Value = 0.2019242199214145
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 35.21571003230599
$ws.Range("R3").Value = 316.941390290754
$ws.Range("S3").Value = 0.0718709816586727
$ws.Range("T3").Value = 0.0718709816586727
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.331241
$ws.Range("H4").Value = 0.9937229999999999
$ws.Range("I4").Value = 0.2019242199214145
$ws.Range("J4").Value = 0.2019242199214145
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 43.50466604598265
$ws.Range("R4").Value = 391.541994413844
$ws.Range("S4").Value = 0.08878773287800031
$ws.Range("T4").Value = 0.08878773287800031
$ws.Range("I5").Value = 0.6029571246591579
$ws.Range("J5").Value = 0.6029571246591579
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 60.37651733792745
$ws.Range("R5").Value = 543.3886560413471
$ws.Range("S5").Value = 0.1232211296102772
$ws.Range("T5").Value = 0.1232211296102772
$ws.Range("I6").Value = 0.6029571246591579
$ws.Range("J6").Value = 0.6029571246591579
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2146108102545086
$ws.Range("T6").Value = 0.2146108102545086
$ws.Range("I7").Value = 0.6029571246591579
$ws.Range("J7").Value = 0.6029571246591579
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 129.9073898047071
$ws.Range("R7").Value = 1169.166508242364
$ws.Range("S7").Value = 0.265125184794372
$ws.Range("T7").Value = 0.265125184794372
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.320077
$ws.Range("H8").Value = 0.9602310000000001
$ws.Range("I8").Value = 0.1951186554194276
$ws.Range("J8").Value = 0.1951186554194276
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 19.53801423035434
$ws.Range("R8").Value = 175.842128073189
$ws.Range("S8").Value = 0.03987471106243465
$ws.Range("T8").Value = 0.03987471106243465
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.320077
$ws.Range("H9").Value = 0.9602310000000001
$ws.Range("I9").Value = 0.1951186554194276
$ws.Range("J9").Value = 0.1951186554194276
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 34.028815333882
$ws.Range("R9").Value = 306.259338004938
$ws.Range("S9").Value = 0.069448673915255
$ws.Range("T9").Value = 0.069448673915255
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.320077
$ws.Range("H10").Value = 0.9602310000000001
$ws.Range("I10").Value = 0.1951186554194276
$ws.Range("J10").Value = 0.1951186554194276
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 42.03840404418533
$ws.Range("R10").Value = 378.345636397668
$ws.Range("S10").Value = 0.08579527044173792
$ws.Range("T10").Value = 0.08579527044173792
